# trading_journal.xlsx update
# Adds exit/closing data for several open trades (rows 8, 17, 22),
# clears the stray "strikethrough" style on some Loss column (L) cells
# (rows 18-20), and fills in the previously-blank trade row 25
# (SBFC Finance) with its full entry/exit data. Also nudges the
# worksheet selection to reflect the reviewed range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Row 8 (Greaves Cotton OAL trade) - trade closed
# ---------------------------------------------------------------
$ws.Range("G8").Value = 365.19
$ws.Range("H8").Value = 53
$ws.Range("L8").Formula = "=(G8-D8)*H8"
$ws.Range("T8").Value = "Closed"

# ---------------------------------------------------------------
# Row 17 (CHOLAHLDNG trade) - date corrected, trade closed
# ---------------------------------------------------------------
$ws.Range("A17").Value = 45841
$ws.Range("G17").Value = 2103.18
$ws.Range("H17").Value = 23
$ws.Range("L17").Formula = "=(G17-D17)*H17"
$ws.Range("T17").Value = "Closed"

# ---------------------------------------------------------------
# Rows 18-20: remove leftover "muted" formatting on the Loss (L)
# column cells - they only ever carried the row's overall style,
# never a value. Re-apply the plain bordered style used elsewhere
# in the L column (e.g. K17) by copying its format across.
# ---------------------------------------------------------------
$ws.Range("K17").Copy()
$ws.Range("L18").PasteSpecial(-4122)
$ws.Range("L19").PasteSpecial(-4122)
$ws.Range("L20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# Row 22 (ENDURANCE trade) - date corrected, trade closed
# ---------------------------------------------------------------
$ws.Range("A22").Value = 45845
$ws.Range("G22").Value = 2667.3
$ws.Range("H22").Value = 11
$ws.Range("L22").Formula = "=(G22-D22)*H22"
$ws.Range("T22").Value = "Closed"

# ---------------------------------------------------------------
# Row 25 (SBFC Finance trade) - previously blank placeholder row,
# now filled in with the full trade record.
# ---------------------------------------------------------------
$ws.Range("A25").Value = 45845
$ws.Range("B25").Value = "SBFC Finance"
$ws.Range("C25").Value = "Long"
$ws.Range("D25").Value = 106.9
$ws.Range("E25").Value = 99.8
$ws.Range("F25").Value = 150
$ws.Range("G25").Value = 113.5
$ws.Range("L25").Formula = "=(G25-D25)*H25"
$ws.Range("M25").Value = "Consodulation Breakout And Retest"
$ws.Range("Q25").Value = "Weekly Day"
$ws.Range("R25").Value = "INR"
$ws.Range("T25").Value = "Closed"

# ---------------------------------------------------------------
# Rows 26-30: blank placeholder rows - extend the Loss (L) formula
# fill-down so they pick up the same (G-D)*H pattern as the rows
# above (they stay #DIV/0! since the rows are otherwise empty).
# ---------------------------------------------------------------
$ws.Range("L26").Formula = "=(G26-D26)*H26"
$ws.Range("L27").Formula = "=(G27-D27)*H27"
$ws.Range("L28").Formula = "=(G28-D28)*H28"
$ws.Range("L29").Formula = "=(G29-D29)*H29"
$ws.Range("L30").Formula = "=(G30-D30)*H30"

# ---------------------------------------------------------------
# Reflect the reviewed range in the sheet selection
# ---------------------------------------------------------------
$ws.Range("A1:U25").Select()
